$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 1 (header row): columns B (2) through DQ (121) hold protein names
# currently stored as the text of a Python 1-tuple repr, e.g. "('ANG_1',)".
# Strip the "('" prefix and "',)" suffix so the cell just contains "ANG_1".
for ($col = 2; $col -le 121; $col++) {
    $cell = $ws.Cells.Item(1, $col)
    $v = $cell.Value2
    if ($v -ne $null -and $v.Length -ge 5 -and $v.Substring(0,2) -eq "('" -and $v.Substring($v.Length-3,3) -eq "',)") {
        $cell.Value = $v.Substring(2, $v.Length - 5)
    }
}

# Column A: rows 2 through 121 hold the same protein names in the same
# tuple-repr form; apply the identical fix-up.
for ($row = 2; $row -le 121; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    $v = $cell.Value2
    if ($v -ne $null -and $v.Length -ge 5 -and $v.Substring(0,2) -eq "('" -and $v.Substring($v.Length-3,3) -eq "',)") {
        $cell.Value = $v.Substring(2, $v.Length - 5)
    }
}
